$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "China"
$ws.Range("B6").Value = 80

$ws.Range("B6").Select()
